# Update "想去人数" (want-to-go count) values across sheets to reflect
# newly generated data (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 5268
$wsExpo.Range("F4").Value = 632
$wsExpo.Range("F5").Value = 305
$wsExpo.Range("F6").Value = 798
$wsExpo.Range("F7").Value = 295
$wsExpo.Range("F8").Value = 11

# Sheet "演出" (Performances)
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 37
$wsShow.Range("F3").Value = 9

# Sheet "全部类型" (All types combined)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 5268
$wsAll.Range("F4").Value = 632
$wsAll.Range("F5").Value = 305
$wsAll.Range("F6").Value = 798
$wsAll.Range("F7").Value = 37
$wsAll.Range("F8").Value = 295
$wsAll.Range("F9").Value = 11
$wsAll.Range("F10").Value = 9

$wb.Save()
